$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 5 de Abril de 2020 a las 22:22'

# Row 4
$ws.Range("B4").Value = 332993
$ws.Range("C4").Value = 21636
$ws.Range("D4").Value = 17018
$ws.Range("E4").Value = 306447
$ws.Range("F4").Value = 8623
$ws.Range("G4").Value = 1076
$ws.Range("H4").Value = 9528

# Row 62 -> Marruecos
$ws.Range("A62").Value = 'Marruecos'
$ws.Range("B62").Value = 1021
$ws.Range("C62").Value = 102
$ws.Range("D62").Value = 76
$ws.Range("E62").Value = 875
$ws.Range("F62").Value = 1
$ws.Range("G62").Value = 11
$ws.Range("H62").Value = 70

# Row 63 -> Eslovenia
$ws.Range("A63").Value = 'Eslovenia'
$ws.Range("B63").Value = 997
$ws.Range("C63").Value = 20
$ws.Range("D63").Value = 79
$ws.Range("E63").Value = 890
$ws.Range("F63").Value = 31
$ws.Range("G63").Value = 6
$ws.Range("H63").Value = 28

# Row 73 -> Camerun
$ws.Range("A73").Value = 'Camerun'
$ws.Range("B73").Value = 650
$ws.Range("C73").Value = 95
$ws.Range("D73").Value = 17
$ws.Range("E73").Value = 624
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 9

# Row 74 -> Azerbaiyan
$ws.Range("A74").Value = 'Azerbaiyan'
$ws.Range("C74").Value = 63
$ws.Range("D74").Value = 32
$ws.Range("E74").Value = 545
$ws.Range("F74").Value = 17
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 7

# Row 75 -> Kazajistan
$ws.Range("A75").Value = 'Kazajistan'
$ws.Range("B75").Value = 584
$ws.Range("C75").Value = 53
$ws.Range("D75").Value = 42
$ws.Range("E75").Value = 536
$ws.Range("F75").Value = 6
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 6

# Row 76 -> Tunez
$ws.Range("A76").Value = 'Tunez'
$ws.Range("B76").Value = 574
$ws.Range("C76").Value = 21
$ws.Range("D76").Value = 5
$ws.Range("E76").Value = 547
$ws.Range("F76").Value = 39
$ws.Range("G76").Value = 4
$ws.Range("H76").Value = 22

# Row 77 -> Bielorrusia
$ws.Range("A77").Value = 'Bielorrusia'
$ws.Range("B77").Value = 562
$ws.Range("C77").Value = 122
$ws.Range("D77").Value = 52
$ws.Range("E77").Value = 502
$ws.Range("F77").Value = 11
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 8

# Row 78 -> Kuwait
$ws.Range("A78").Value = 'Kuwait'
$ws.Range("B78").Value = 556
$ws.Range("C78").Value = 77
$ws.Range("D78").Value = 99
$ws.Range("E78").Value = 456
$ws.Range("F78").Value = 17
$ws.Range("H78").Value = 1

# Row 99
$ws.Range("B99").Value = 261
$ws.Range("C99").Value = 16
$ws.Range("D99").Value = 37
$ws.Range("E99").Value = 221
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 3

# Row 108 -> Niger
$ws.Range("A108").Value = 'Niger'
$ws.Range("B108").Value = 184
$ws.Range("C108").Value = 40
$ws.Range("D108").Value = 13
$ws.Range("E108").Value = 161
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 2
$ws.Range("H108").Value = 10

# Row 109 -> Islas Feroe
$ws.Range("A109").Value = 'Islas Feroe'
$ws.Range("B109").Value = 181
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 99
$ws.Range("E109").Value = 82
$ws.Range("F109").Value = 1
$ws.Range("H109").Value = 0

# Row 110 -> Sri Lanka
$ws.Range("A110").Value = 'Sri Lanka'
$ws.Range("B110").Value = 176
$ws.Range("C110").Value = 10
$ws.Range("D110").Value = 33
$ws.Range("E110").Value = 138
$ws.Range("F110").Value = 5
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 5

# Row 111 -> Georgia
$ws.Range("A111").Value = 'Georgia'
$ws.Range("B111").Value = 174
$ws.Range("C111").Value = 12
$ws.Range("D111").Value = 36
$ws.Range("E111").Value = 136
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 2

# Row 112 -> Venezuela
$ws.Range("A112").Value = 'Venezuela'
$ws.Range("B112").Value = 159
$ws.Range("C112").Value = 4
$ws.Range("D112").Value = 52
$ws.Range("E112").Value = 100
$ws.Range("F112").Value = 6
$ws.Range("H112").Value = 7

# Row 113 -> Bolivia
$ws.Range("A113").Value = 'Bolivia'
$ws.Range("B113").Value = 157
$ws.Range("C113").Value = 18
$ws.Range("D113").Value = 2
$ws.Range("E113").Value = 145
$ws.Range("F113").Value = 3
$ws.Range("H113").Value = 10

# Row 114 -> Consejo Danes para los Refugiados
$ws.Range("A114").Value = 'Consejo Danes para los Refugiados'
$ws.Range("B114").Value = 154
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 3
$ws.Range("E114").Value = 133
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 18

# Row 115 -> Martinica
$ws.Range("A115").Value = 'Martinica'
$ws.Range("B115").Value = 149
$ws.Range("C115").Value = 4
$ws.Range("D115").Value = 50
$ws.Range("E115").Value = 95
$ws.Range("F115").Value = 21
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 4

# Row 116 -> Kirguistan
$ws.Range("A116").Value = 'Kirguistan'
$ws.Range("B116").Value = 147
$ws.Range("C116").Value = 3
$ws.Range("D116").Value = 9
$ws.Range("E116").Value = 137
$ws.Range("F116").Value = 5
$ws.Range("H116").Value = 1

# Row 154
$ws.Range("D154").Value = 1
$ws.Range("E154").Value = 20

# Row 174 -> Granada
$ws.Range("A174").Value = 'Granada'
$ws.Range("F174").Value = 2

# Row 175 -> Fiyi
$ws.Range("A175").Value = 'Fiyi'
$ws.Range("F175").Value = 0

# Row 182 -> Mozambique
$ws.Range("A182").Value = 'Mozambique'
$ws.Range("D182").Value = 1
$ws.Range("H182").Value = 0

# Row 183 -> Surinam
$ws.Range("A183").Value = 'Surinam'
$ws.Range("D183").Value = 0
$ws.Range("H183").Value = 1

# Row 192 -> Somalia
$ws.Range("A192").Value = 'Somalia'

# Row 193 -> San Vicente y las Granadinas
$ws.Range("A193").Value = 'San Vicente y las Granadinas'

# Row 196 -> San Bartolome
$ws.Range("A196").Value = 'San Bartolome'
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 1
$ws.Range("H196").Value = 0

# Row 197 -> Botsuana
$ws.Range("A197").Value = 'Botsuana'
$ws.Range("C197").Value = 2
$ws.Range("D197").Value = 0
$ws.Range("H197").Value = 1

# Row 207 -> Burundi
$ws.Range("A207").Value = 'Burundi'

# Row 208 -> Anguila
$ws.Range("A208").Value = 'Anguila'

# Row 209 -> Bonaire, San Eustaquio y Saba
$ws.Range("A209").Value = 'Bonaire, San Eustaquio y Saba'
$ws.Range("C209").Value = 0

# Row 210 -> Islas Malvinas
$ws.Range("A210").Value = 'Islas Malvinas'
$ws.Range("C210").Value = 1
